# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1) Update the time_taken (column F) timestamps on the "data" sheet ---
$ws.Range("F2").Value  = "2021-10-05 14:21:51.970787"
$ws.Range("F3").Value  = "2021-10-05 14:21:51.970796"
$ws.Range("F4").Value  = "2021-10-05 14:21:51.970800"
$ws.Range("F5").Value  = "2021-10-05 14:21:51.970803"
$ws.Range("F6").Value  = "2021-10-05 14:21:51.970805"
$ws.Range("F7").Value  = "2021-10-05 14:21:51.970808"
$ws.Range("F8").Value  = "2021-10-05 14:21:51.970811"
$ws.Range("F9").Value  = "2021-10-05 14:21:51.970814"
$ws.Range("F10").Value = "2021-10-05 14:21:51.970817"
$ws.Range("F11").Value = "2021-10-05 14:21:51.970819"
$ws.Range("F12").Value = "2021-10-05 14:21:51.970822"
$ws.Range("F13").Value = "2021-10-05 14:21:51.970825"
$ws.Range("F14").Value = "2021-10-05 14:21:51.970827"

# --- 2) Add a new "metadata" worksheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the bold/bordered header style from the "data" sheet header onto B1:G1
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row 2
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Non-syndromic hypotrichosis"
$meta.Range("C2").Value = 189

# data_version must be stored as text "1.9" (not the number 1.9).
# Start from the plain/unstyled format used by the "data" sheet's text cells,
# then force text so the numeric-looking string isn't parsed as a number,
# and finally restore the plain format so no stray number format lingers.
$ws.Range("B2").Copy()
$meta.Range("D2").PasteSpecial(-4122)
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.9"
$ws.Range("B2").Copy()
$meta.Range("D2").PasteSpecial(-4122)

$meta.Range("E2").Value = "2021-07-28T14:09:42.959731Z"
$meta.Range("F2").Value = "2021-10-05 14:21:51.967575"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/189/?format=json"

# Apply the bold/bordered style (same as the header row) to A2, matching the "data" sheet's A-column style
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$excel.CutCopyMode = $false
